$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '24.046.96'
Set-TextValue $ws.Range('E2') '  -3.66%  '
Set-TextValue $ws.Range('D3') '1.637.01'
Set-TextValue $ws.Range('E3') '  -3.34%  '
Set-TextValue $ws.Range('E4') '  -0.21%  '
Set-TextValue $ws.Range('D5') '1.000'
Set-TextValue $ws.Range('E5') '  -0.49%  '
Set-TextValue $ws.Range('D6') '305.91'
Set-TextValue $ws.Range('E6') '  -3.46%  '
Set-TextValue $ws.Range('D7') '0.3866'
Set-TextValue $ws.Range('E7') '  -2.38%  '
Set-TextValue $ws.Range('E8') '  -4.22%  '
Set-TextValue $ws.Range('D9') '0.9998'
Set-TextValue $ws.Range('E9') '  -0.34%  '
Set-TextValue $ws.Range('E10') '  -6.53%  '
Set-TextValue $ws.Range('D11') '1.339'
Set-TextValue $ws.Range('E11') '  -7.35%  '
Set-TextValue $ws.Range('D12') '0.08455'
Set-TextValue $ws.Range('E12') '  -3.26%  '
Set-TextValue $ws.Range('D13') '23.53'
Set-TextValue $ws.Range('E13') '  -7.96%  '
Set-TextValue $ws.Range('E14') '  -4.77%  '
Set-TextValue $ws.Range('E15') '  -5.18%  '
Set-TextValue $ws.Range('D16') '7.417'
Set-TextValue $ws.Range('E16') '  -5.84%  '
Set-TextValue $ws.Range('D17') '1.638.02'
Set-TextValue $ws.Range('E17') '  -16.07%  '
Set-TextValue $ws.Range('D18') '94.84'
Set-TextValue $ws.Range('E18') '  -0.03%  '
Set-TextValue $ws.Range('D19') '0.06866'
Set-TextValue $ws.Range('E19') '  -5.19%  '
Set-TextValue $ws.Range('D20') '20.60'
Set-TextValue $ws.Range('E20') '  +0.67%  '
Set-TextValue $ws.Range('D21') '6.865'
Set-TextValue $ws.Range('E21') '  -4.45%  '
Set-TextValue $ws.Range('D22') '1.000'
Set-TextValue $ws.Range('E22') '  -0.53%  '
Set-TextValue $ws.Range('E23') '  -5.19%  '
Set-TextValue $ws.Range('D24') '24.086.12'
Set-TextValue $ws.Range('E24') '  -3.40%  '
Set-TextValue $ws.Range('E25') '  -3.25%  '
Set-TextValue $ws.Range('D26') '2.672'
Set-TextValue $ws.Range('E26') '  -6.45%  '
Set-TextValue $ws.Range('D27') '22.25'
Set-TextValue $ws.Range('E27') '  -3.84%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D28') '156.76'
Set-TextValue $ws.Range('E28') '  -3.29%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D29') '8.654'
Set-TextValue $ws.Range('E29') '  +6.79%  '
Set-TextValue $ws.Range('D30') '139.61'
Set-TextValue $ws.Range('E30') '  -6.14%  '
Set-TextValue $ws.Range('D31') '5.331'
Set-TextValue $ws.Range('E31') '  -11.90%  '
Set-TextValue $ws.Range('D32') '2.416'
Set-TextValue $ws.Range('E32') '  -7.55%  '
Set-TextValue $ws.Range('D33') '1.812.31'
Set-TextValue $ws.Range('E33') '  -16.02%  '
Set-TextValue $ws.Range('D34') '6.886'
Set-TextValue $ws.Range('E34') '  -2.84%  '
Set-TextValue $ws.Range('E35') '  -6.63%  '
Set-TextValue $ws.Range('D36') '0.02868'
Set-TextValue $ws.Range('E36') '  -7.92%  '
Set-TextValue $ws.Range('D37') '0.2665'
Set-TextValue $ws.Range('E37') '  -7.10%  '
Set-TextValue $ws.Range('D38') '0.9459'
Set-TextValue $ws.Range('E38') '  -8.83%  '
Set-TextValue $ws.Range('D39') '0.09140'
Set-TextValue $ws.Range('E39') '  -5.47%  '
Set-TextValue $ws.Range('D40') '1.445'
Set-TextValue $ws.Range('E40') '  -1.96%  '
Set-TextValue $ws.Range('D41') '9.835'
Set-TextValue $ws.Range('E41') '  -9.26%  '
Set-TextValue $ws.Range('D42') '0.7496'
Set-TextValue $ws.Range('E42') '  -7.54%  '
Set-TextValue $ws.Range('D43') '12.93'
Set-TextValue $ws.Range('E43') '  -7.16%  '
Set-TextValue $ws.Range('D44') '15.90'
Set-TextValue $ws.Range('E44') '  -5.78%  '
Set-TextValue $ws.Range('D45') '0.6840'
Set-TextValue $ws.Range('E45') '  -6.12%  '
Set-TextValue $ws.Range('D46') '2.449'
Set-TextValue $ws.Range('E46') '  -7.16%  '
Set-TextValue $ws.Range('D47') '4.081'
Set-TextValue $ws.Range('E47') '  -3.44%  '
Set-TextValue $ws.Range('D48') '0.9998'
Set-TextValue $ws.Range('E48') '  -0.29%  '
Set-TextValue $ws.Range('D49') '0.08344'
Set-TextValue $ws.Range('E49') '  -6.76%  '
Set-TextValue $ws.Range('E50') '  -10.11%  '
Set-TextValue $ws.Range('D51') '131.99'
Set-TextValue $ws.Range('E51') '  -5.78%  '
